# Update the "Generate Report for Handback" timestamps that were refreshed
# by the latest handback run.
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 3d6c9228... row.
$wsOverview.Range("G3").Value = "2016-08-16 04:43:57"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the 3d6c9228... row.
$wsZhCn.Range("H3").Value = "2016-08-16 04:43:52"
$wsZhCn.Range("K3").Value = "2016-08-16 04:44:15"

# de-de sheet: "Correspond Handoff Datetime" (matches the Overview value
# above) and "Correspond Handback DateTime" for the 3d6c9228... row.
$wsDeDe.Range("H3").Value = "2016-08-16 04:43:57"
$wsDeDe.Range("K3").Value = "2016-08-16 04:44:22"
